$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper now also captures player "height" and "weight", inserted as
# new columns E and F. The pre-existing "fantasy points" column (old E)
# shifts right to become column G.

# Save the existing "fantasy points" values (column E) before they are
# overwritten below.
$fantasyPoints = @{}
for ($r = 2; $r -le 13; $r++) {
    $fantasyPoints[$r] = $ws.Cells.Item($r, 5).Value2
}

# Header row: E1 becomes "height", F1 (new) becomes "weight", G1 (new)
# becomes "fantasy points" (moved from the old E1).
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# F1/G1 are brand-new header cells - copy E1's header formatting (bold,
# border, centered) onto them so they match B1:E1.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Data rows: height is constant (6'4" => 6.333333333333333 ft) and weight
# is constant (245 lbs) for every row; fantasy points move from E to G.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 5).Value2 = 6.333333333333333
    $ws.Cells.Item($r, 6).Value2 = 245
    $ws.Cells.Item($r, 7).Value2 = $fantasyPoints[$r]
}
